$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.18897966666667
$ws.Range("H2").Value = 45.566939
$ws.Range("I2").Value = 0.01327029680642083
$ws.Range("J2").Value = 0.01327029680642082
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8690176666666667
$ws.Range("N2").Value = 2.607053
$ws.Range("O2").Value = 0.004684547990100867
$ws.Range("P2").Value = 0.004684547990100866
$ws.Range("Q2").Value = 13.19949166897411
$ws.Range("R2").Value = 118.795425020767
$ws.Range("S2").Value = 0.00006216534223256062
$ws.Range("T2").Value = 0.0000621653422325606

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.18897966666667
$ws.Range("H3").Value = 45.566939
$ws.Range("I3").Value = 0.01327029680642083
$ws.Range("J3").Value = 0.01327029680642082
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.730875666666667
$ws.Range("N3").Value = 23.192627
$ws.Range("O3").Value = 0.04167424835552215
$ws.Range("P3").Value = 0.04167424835552214
$ws.Range("Q3").Value = 117.4241133065281
$ws.Range("R3").Value = 1056.817019758753
$ws.Range("S3").Value = 0.000553029644862274
$ws.Range("T3").Value = 0.0005530296448622736

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.18897966666667
$ws.Range("H4").Value = 45.566939
$ws.Range("I4").Value = 0.01327029680642083
$ws.Range("J4").Value = 0.01327029680642082
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 67.01802199999999
$ws.Range("N4").Value = 201.054066
$ws.Range("O4").Value = 0.3612689963655924
$ws.Range("P4").Value = 0.3612689963655923
$ws.Range("Q4").Value = 1017.935373458219
$ws.Range("R4").Value = 9161.418361123973
$ws.Range("S4").Value = 0.004794146808729177
$ws.Range("T4").Value = 0.004794146808729176

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.18897966666667
$ws.Range("H5").Value = 45.566939
$ws.Range("I5").Value = 0.01327029680642083
$ws.Range("J5").Value = 0.01327029680642082
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 79.20011933333333
$ws.Range("N5").Value = 237.600358
$ws.Range("O5").Value = 0.4269381096265194
$ws.Range("P5").Value = 0.4269381096265193
$ws.Range("Q5").Value = 1202.969002151573
$ws.Range("R5").Value = 10826.72101936416
$ws.Range("S5").Value = 0.005665595432716145
$ws.Range("T5").Value = 0.005665595432716141

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.18897966666667
$ws.Range("H6").Value = 45.566939
$ws.Range("I6").Value = 0.01327029680642083
$ws.Range("J6").Value = 0.01327029680642082
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.17078333333333
$ws.Range("N6").Value = 81.51235
$ws.Range("O6").Value = 0.1464674923604922
$ws.Range("P6").Value = 0.1464674923604922
$ws.Range("Q6").Value = 412.6964755774055
$ws.Range("R6").Value = 3714.26828019665
$ws.Range("S6").Value = 0.001943667096115907
$ws.Range("T6").Value = 0.001943667096115906

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.18897966666667
$ws.Range("H7").Value = 45.566939
$ws.Range("I7").Value = 0.01327029680642083
$ws.Range("J7").Value = 0.01327029680642082
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.518443
$ws.Range("N7").Value = 10.555329
$ws.Range("O7").Value = 0.01896660530177307
$ws.Range("P7").Value = 0.01896660530177307
$ws.Range("Q7").Value = 53.44155918532567
$ws.Range("R7").Value = 480.974032667931
$ws.Range("S7").Value = 0.0002516924817647635
$ws.Range("T7").Value = 0.0002516924817647634

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 95.50314333333334
$ws.Range("H8").Value = 286.50943
$ws.Range("I8").Value = 0.08343911742543099
$ws.Range("J8").Value = 0.08343911742543098
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8690176666666667
$ws.Range("N8").Value = 2.607053
$ws.Range("O8").Value = 0.004684547990100867
$ws.Range("P8").Value = 0.004684547990100866
$ws.Range("Q8").Value = 82.99391877886556
$ws.Range("R8").Value = 746.9452690097901
$ws.Range("S8").Value = 0.0003908745498310929
$ws.Range("T8").Value = 0.0003908745498310928

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 95.50314333333334
$ws.Range("H9").Value = 286.50943
$ws.Range("I9").Value = 0.08343911742543099
$ws.Range("J9").Value = 0.08343911742543098
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.730875666666667
$ws.Range("N9").Value = 23.192627
$ws.Range("O9").Value = 0.04167424835552215
$ws.Range("P9").Value = 0.04167424835552214
$ws.Range("Q9").Value = 738.3229268858456
$ws.Range("R9").Value = 6644.906341972611
$ws.Range("S9").Value = 0.003477262502152987
$ws.Range("T9").Value = 0.003477262502152986

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 95.50314333333334
$ws.Range("H10").Value = 286.50943
$ws.Range("I10").Value = 0.08343911742543099
$ws.Range("J10").Value = 0.08343911742543098
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 67.01802199999999
$ws.Range("N10").Value = 201.054066
$ws.Range("O10").Value = 0.3612689963655924
$ws.Range("P10").Value = 0.3612689963655923
$ws.Range("Q10").Value = 6400.431760982486
$ws.Range("R10").Value = 57603.88584884237
$ws.Range("S10").Value = 0.03014396620991627
$ws.Range("T10").Value = 0.03014396620991626

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 95.50314333333334
$ws.Range("H11").Value = 286.50943
$ws.Range("I11").Value = 0.08343911742543099
$ws.Range("J11").Value = 0.08343911742543098
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 79.20011933333333
$ws.Range("N11").Value = 237.600358
$ws.Range("O11").Value = 0.4269381096265194
$ws.Range("P11").Value = 0.4269381096265193
$ws.Range("Q11").Value = 7563.860348708438
$ws.Range("R11").Value = 68074.74313837594
$ws.Range("S11").Value = 0.03562333906251868
$ws.Range("T11").Value = 0.03562333906251867

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 95.50314333333334
$ws.Range("H12").Value = 286.50943
$ws.Range("I12").Value = 0.08343911742543099
$ws.Range("J12").Value = 0.08343911742543098
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 27.17078333333333
$ws.Range("N12").Value = 81.51235
$ws.Range("O12").Value = 0.1464674923604922
$ws.Range("P12").Value = 0.1464674923604922
$ws.Range("Q12").Value = 2594.895215162278
$ws.Range("R12").Value = 23354.0569364605
$ws.Range("S12").Value = 0.01222111829407553
$ws.Range("T12").Value = 0.01222111829407552

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 95.50314333333334
$ws.Range("H13").Value = 286.50943
$ws.Range("I13").Value = 0.08343911742543099
$ws.Range("J13").Value = 0.08343911742543098
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.518443
$ws.Range("N13").Value = 10.555329
$ws.Range("O13").Value = 0.01896660530177307
$ws.Range("P13").Value = 0.01896660530177307
$ws.Range("Q13").Value = 336.0223661391634
$ws.Range("R13").Value = 3024.20129525247
$ws.Range("S13").Value = 0.001582556806936445
$ws.Range("T13").Value = 0.001582556806936444

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 402.0913113333333
$ws.Range("H14").Value = 1206.273934
$ws.Range("I14").Value = 0.3512988470441011
$ws.Range("J14").Value = 0.351298847044101
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.8690176666666667
$ws.Range("N14").Value = 2.607053
$ws.Range("O14").Value = 0.004684547990100867
$ws.Range("P14").Value = 0.004684547990100866
$ws.Range("Q14").Value = 349.4244531618335
$ws.Range("R14").Value = 3144.820078456502
$ws.Range("S14").Value = 0.001645676307845196
$ws.Range("T14").Value = 0.001645676307845195

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 402.0913113333333
$ws.Range("H15").Value = 1206.273934
$ws.Range("I15").Value = 0.3512988470441011
$ws.Range("J15").Value = 0.351298847044101
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.730875666666667
$ws.Range("N15").Value = 23.192627
$ws.Range("O15").Value = 0.04167424835552215
$ws.Range("P15").Value = 0.04167424835552214
$ws.Range("Q15").Value = 3108.517934564957
$ws.Range("R15").Value = 27976.66141108462
$ws.Range("S15").Value = 0.01464011539872446
$ws.Range("T15").Value = 0.01464011539872445

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 402.0913113333333
$ws.Range("H16").Value = 1206.273934
$ws.Range("I16").Value = 0.3512988470441011
$ws.Range("J16").Value = 0.351298847044101
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 67.01802199999999
$ws.Range("N16").Value = 201.054066
$ws.Range("O16").Value = 0.3612689963655924
$ws.Range("P16").Value = 0.3612689963655923
$ws.Range("Q16").Value = 26947.36434894618
$ws.Range("R16").Value = 242526.2791405156
$ws.Range("S16").Value = 0.1269133818960121
$ws.Range("T16").Value = 0.1269133818960121

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 402.0913113333333
$ws.Range("H17").Value = 1206.273934
$ws.Range("I17").Value = 0.3512988470441011
$ws.Range("J17").Value = 0.351298847044101
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 79.20011933333333
$ws.Range("N17").Value = 237.600358
$ws.Range("O17").Value = 0.4269381096265194
$ws.Range("P17").Value = 0.4269381096265193
$ws.Range("Q17").Value = 31845.67984049648
$ws.Range("R17").Value = 286611.1185644683
$ws.Range("S17").Value = 0.1499828656709843
$ws.Range("T17").Value = 0.1499828656709842

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 402.0913113333333
$ws.Range("H18").Value = 1206.273934
$ws.Range("I18").Value = 0.3512988470441011
$ws.Range("J18").Value = 0.351298847044101
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 27.17078333333333
$ws.Range("N18").Value = 81.51235
$ws.Range("O18").Value = 0.1464674923604922
$ws.Range("P18").Value = 0.1464674923604922
$ws.Range("Q18").Value = 10925.13590045388
$ws.Range("R18").Value = 98326.22310408489
$ws.Range("S18").Value = 0.0514538611956816
$ws.Range("T18").Value = 0.05145386119568158

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 402.0913113333333
$ws.Range("H19").Value = 1206.273934
$ws.Range("I19").Value = 0.3512988470441011
$ws.Range("J19").Value = 0.351298847044101
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 3.518443
$ws.Range("N19").Value = 10.555329
$ws.Range("O19").Value = 0.01896660530177307
$ws.Range("P19").Value = 0.01896660530177307
$ws.Range("Q19").Value = 1414.735359721587
$ws.Range("R19").Value = 12732.61823749429
$ws.Range("S19").Value = 0.006662946574853414
$ws.Range("T19").Value = 0.006662946574853411

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 192.2123056666667
$ws.Range("H20").Value = 576.636917
$ws.Range("I20").Value = 0.1679319086614409
$ws.Range("J20").Value = 0.1679319086614409
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 0.8690176666666667
$ws.Range("N20").Value = 2.607053
$ws.Range("O20").Value = 0.004684547990100867
$ws.Range("P20").Value = 0.004684547990100866
$ws.Range("Q20").Value = 167.0358893750668
$ws.Range("R20").Value = 1503.323004375601
$ws.Range("S20").Value = 0.0007866850851937553
$ws.Range("T20").Value = 0.0007866850851937552

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 192.2123056666667
$ws.Range("H21").Value = 576.636917
$ws.Range("I21").Value = 0.1679319086614409
$ws.Range("J21").Value = 0.1679319086614409
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 7.730875666666667
$ws.Range("N21").Value = 23.192627
$ws.Range("O21").Value = 0.04167424835552215
$ws.Range("P21").Value = 0.04167424835552214
$ws.Range("Q21").Value = 1485.969436712329
$ws.Range("R21").Value = 13373.72493041096
$ws.Range("S21").Value = 0.00699843606837375
$ws.Range("T21").Value = 0.006998436068373749

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 192.2123056666667
$ws.Range("H22").Value = 576.636917
$ws.Range("I22").Value = 0.1679319086614409
$ws.Range("J22").Value = 0.1679319086614409
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 67.01802199999999
$ws.Range("N22").Value = 201.054066
$ws.Range("O22").Value = 0.3612689963655924
$ws.Range("P22").Value = 0.3612689963655923
$ws.Range("Q22").Value = 12881.68852983939
$ws.Range("R22").Value = 115935.1967685545
$ws.Range("S22").Value = 0.06066859209987709
$ws.Range("T22").Value = 0.06066859209987709

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 192.2123056666667
$ws.Range("H23").Value = 576.636917
$ws.Range("I23").Value = 0.1679319086614409
$ws.Range("J23").Value = 0.1679319086614409
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 79.20011933333333
$ws.Range("N23").Value = 237.600358
$ws.Range("O23").Value = 0.4269381096265194
$ws.Range("P23").Value = 0.4269381096265193
$ws.Range("Q23").Value = 15223.23754613514
$ws.Range("R23").Value = 137009.1379152163
$ws.Range("S23").Value = 0.07169653162988891
$ws.Range("T23").Value = 0.07169653162988888

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 192.2123056666667
$ws.Range("H24").Value = 576.636917
$ws.Range("I24").Value = 0.1679319086614409
$ws.Range("J24").Value = 0.1679319086614409
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 27.17078333333333
$ws.Range("N24").Value = 81.51235
$ws.Range("O24").Value = 0.1464674923604922
$ws.Range("P24").Value = 0.1464674923604922
$ws.Range("Q24").Value = 5222.558911269439
$ws.Range("R24").Value = 47003.03020142495
$ws.Range("S24").Value = 0.02459656554895247
$ws.Range("T24").Value = 0.02459656554895247

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 192.2123056666667
$ws.Range("H25").Value = 576.636917
$ws.Range("I25").Value = 0.1679319086614409
$ws.Range("J25").Value = 0.1679319086614409
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 3.518443
$ws.Range("N25").Value = 10.555329
$ws.Range("O25").Value = 0.01896660530177307
$ws.Range("P25").Value = 0.01896660530177307
$ws.Range("Q25").Value = 676.2880413867437
$ws.Range("R25").Value = 6086.592372480693
$ws.Range("S25").Value = 0.003185098229154956
$ws.Range("T25").Value = 0.003185098229154956

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 402.444082
$ws.Range("H26").Value = 1207.332246
$ws.Range("I26").Value = 0.3516070554658648
$ws.Range("J26").Value = 0.3516070554658648
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 0.8690176666666667
$ws.Range("N26").Value = 2.607053
$ws.Range("O26").Value = 0.004684547990100867
$ws.Range("P26").Value = 0.004684547990100866
$ws.Range("Q26").Value = 349.7310171034487
$ws.Range("R26").Value = 3147.579153931038
$ws.Range("S26").Value = 0.001647120124987901
$ws.Range("T26").Value = 0.001647120124987901

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 402.444082
$ws.Range("H27").Value = 1207.332246
$ws.Range("I27").Value = 0.3516070554658648
$ws.Range("J27").Value = 0.3516070554658648
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 7.730875666666667
$ws.Range("N27").Value = 23.192627
$ws.Range("O27").Value = 0.04167424835552215
$ws.Range("P27").Value = 0.04167424835552214
$ws.Range("Q27").Value = 3111.245160727804
$ws.Range("R27").Value = 28001.20644655024
$ws.Range("S27").Value = 0.0146529597530383
$ws.Range("T27").Value = 0.0146529597530383

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 402.444082
$ws.Range("H28").Value = 1207.332246
$ws.Range("I28").Value = 0.3516070554658648
$ws.Range("J28").Value = 0.3516070554658648
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 67.01802199999999
$ws.Range("N28").Value = 201.054066
$ws.Range("O28").Value = 0.3612689963655924
$ws.Range("P28").Value = 0.3612689963655923
$ws.Range("Q28").Value = 26971.0063412458
$ws.Range("R28").Value = 242739.0570712122
$ws.Range("S28").Value = 0.1270247280432142
$ws.Range("T28").Value = 0.1270247280432141

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 402.444082
$ws.Range("H29").Value = 1207.332246
$ws.Range("I29").Value = 0.3516070554658648
$ws.Range("J29").Value = 0.3516070554658648
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 79.20011933333333
$ws.Range("N29").Value = 237.600358
$ws.Range("O29").Value = 0.4269381096265194
$ws.Range("P29").Value = 0.4269381096265193
$ws.Range("Q29").Value = 31873.61931939378
$ws.Range("R29").Value = 286862.5738745441
$ws.Range("S29").Value = 0.1501144515919431
$ws.Range("T29").Value = 0.150114451591943

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 402.444082
$ws.Range("H30").Value = 1207.332246
$ws.Range("I30").Value = 0.3516070554658648
$ws.Range("J30").Value = 0.3516070554658648
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 27.17078333333333
$ws.Range("N30").Value = 81.51235
$ws.Range("O30").Value = 0.1464674923604922
$ws.Range("P30").Value = 0.1464674923604922
$ws.Range("Q30").Value = 10934.72095580423
$ws.Range("R30").Value = 98412.4886022381
$ws.Range("S30").Value = 0.05149900371034172
$ws.Range("T30").Value = 0.05149900371034171

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 402.444082
$ws.Range("H31").Value = 1207.332246
$ws.Range("I31").Value = 0.3516070554658648
$ws.Range("J31").Value = 0.3516070554658648
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 3.518443
$ws.Range("N31").Value = 10.555329
$ws.Range("O31").Value = 0.01896660530177307
$ws.Range("P31").Value = 0.01896660530177307
$ws.Range("Q31").Value = 1415.976563204326
$ws.Range("R31").Value = 12743.78906883893
$ws.Range("S31").Value = 0.00666879224233969
$ws.Range("T31").Value = 0.006668792242339688

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 37.14495166666666
$ws.Range("H32").Value = 111.434855
$ws.Range("I32").Value = 0.03245277459674146
$ws.Range("J32").Value = 0.03245277459674145
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 0.8690176666666667
$ws.Range("N32").Value = 2.607053
$ws.Range("O32").Value = 0.004684547990100867
$ws.Range("P32").Value = 0.004684547990100866
$ws.Range("Q32").Value = 32.27961922581277
$ws.Range("R32").Value = 290.516573032315
$ws.Range("S32").Value = 0.0001520265800103617
$ws.Range("T32").Value = 0.0001520265800103616

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 37.14495166666666
$ws.Range("H33").Value = 111.434855
$ws.Range("I33").Value = 0.03245277459674146
$ws.Range("J33").Value = 0.03245277459674145
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 7.730875666666667
$ws.Range("N33").Value = 23.192627
$ws.Range("O33").Value = 0.04167424835552215
$ws.Range("P33").Value = 0.04167424835552214
$ws.Range("Q33").Value = 287.1630029793428
$ws.Range("R33").Value = 2584.467026814085
$ws.Range("S33").Value = 0.001352444988370384
$ws.Range("T33").Value = 0.001352444988370383

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 37.14495166666666
$ws.Range("H34").Value = 111.434855
$ws.Range("I34").Value = 0.03245277459674146
$ws.Range("J34").Value = 0.03245277459674145
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 67.01802199999999
$ws.Range("N34").Value = 201.054066
$ws.Range("O34").Value = 0.3612689963655924
$ws.Range("P34").Value = 0.3612689963655923
$ws.Range("Q34").Value = 2489.381187985603
$ws.Range("R34").Value = 22404.43069187043
$ws.Range("S34").Value = 0.01172418130784358
$ws.Range("T34").Value = 0.01172418130784358

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 37.14495166666666
$ws.Range("H35").Value = 111.434855
$ws.Range("I35").Value = 0.03245277459674146
$ws.Range("J35").Value = 0.03245277459674145
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 79.20011933333333
$ws.Range("N35").Value = 237.600358
$ws.Range("O35").Value = 0.4269381096265194
$ws.Range("P35").Value = 0.4269381096265193
$ws.Range("Q35").Value = 2941.884604630899
$ws.Range("R35").Value = 26476.96144167809
$ws.Range("S35").Value = 0.01385532623846833
$ws.Range("T35").Value = 0.01385532623846832

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 37.14495166666666
$ws.Range("H36").Value = 111.434855
$ws.Range("I36").Value = 0.03245277459674146
$ws.Range("J36").Value = 0.03245277459674145
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 27.17078333333333
$ws.Range("N36").Value = 81.51235
$ws.Range("O36").Value = 0.1464674923604922
$ws.Range("P36").Value = 0.1464674923604922
$ws.Range("Q36").Value = 1009.257433662139
$ws.Range("R36").Value = 9083.31690295925
$ws.Range("S36").Value = 0.004753276515325006
$ws.Range("T36").Value = 0.004753276515325004

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 37.14495166666666
$ws.Range("H37").Value = 111.434855
$ws.Range("I37").Value = 0.03245277459674146
$ws.Range("J37").Value = 0.03245277459674145
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 3.518443
$ws.Range("N37").Value = 10.555329
$ws.Range("O37").Value = 0.01896660530177307
$ws.Range("P37").Value = 0.01896660530177307
$ws.Range("Q37").Value = 130.6923951769217
$ws.Range("R37").Value = 1176.231556592295
$ws.Range("S37").Value = 0.0006155189667238029
$ws.Range("T37").Value = 0.0006155189667238027

